# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for rows 2-50 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1;
    3 = 1;
    4 = 2;
    5 = 1;
    6 = 1;
    7 = 0;
    8 = 0;
    9 = 2;
    10 = 2;
    11 = 0;
    12 = 1;
    13 = 0;
    14 = 1;
    15 = 1;
    16 = 1;
    17 = 4;
    18 = 0;
    19 = 2;
    20 = 1;
    21 = 0;
    22 = 1;
    23 = 0;
    24 = 2;
    25 = 2;
    26 = 1;
    27 = 1;
    28 = 1;
    29 = 1;
    30 = 1;
    31 = 2;
    32 = 1;
    33 = 0;
    34 = 1;
    35 = 2;
    36 = 1;
    37 = 1;
    38 = 1;
    39 = 1;
    40 = 0;
    41 = 1;
    42 = 0;
    43 = 0;
    44 = 2;
    45 = 1;
    46 = 1;
    47 = 1;
    48 = 3;
    49 = 2;
    50 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
